$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-25 Thursday" "2025-12-26 Friday"

Replace-Text "731÷5=" "909÷8="
Replace-Text "301÷7=" "389÷5="
Replace-Text "487÷2=" "337÷3="
Replace-Text "778÷7=" "348÷8="
Replace-Text "948÷6=" "825÷4="
Replace-Text "289÷4=" "476÷3="
Replace-Text "490÷8=" "108÷8="
Replace-Text "564÷4=" "119÷9="
Replace-Text "550÷7=" "487÷5="
Replace-Text "164÷9=" "175÷8="
Replace-Text "872÷2=" "483÷9="
Replace-Text "151÷9=" "745÷4="
Replace-Text "165÷6=" "885÷4="
Replace-Text "225÷3=" "278÷3="
Replace-Text "508÷5=" "134÷2="
Replace-Text "538÷3=" "137÷5="
Replace-Text "668÷4=" "316÷2="
Replace-Text "507÷6=" "192÷5="
Replace-Text "674÷5=" "682÷3="
Replace-Text "222÷7=" "344÷3="
Replace-Text "491÷4=" "737÷4="
Replace-Text "330÷4=" "186÷3="
Replace-Text "722÷9=" "215÷8="
Replace-Text "157÷9=" "675÷7="
Replace-Text "992÷4=" "820÷3="
